$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.374.17"
$ws.Range("E2").Value = "  -2.99%  "
$ws.Range("D3").Value = "1.981.66"
$ws.Range("E3").Value = "  -3.54%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "243.66"
$ws.Range("E5").Value = "  -3.63%  "
$ws.Range("E6").Value = "  -3.62%  "
$ws.Range("D7").Value = "58.82"
$ws.Range("E7").Value = "  -11.62%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "57.45"
$ws.Range("E10").Value = "  -3.98%  "
$ws.Range("D11").Value = "0.0817"
$ws.Range("E11").Value = "  +6.36%  "
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "23.79"
$ws.Range("E13").Value = "  +8.03%  "
$ws.Range("D14").Value = "0.861"
$ws.Range("E14").Value = "  -5.46%  "
$ws.Range("D15").Value = "14.02"
$ws.Range("E15").Value = "  -6.25%  "
$ws.Range("D16").Value = "2.271.32"
$ws.Range("E16").Value = "  -3.53%  "
$ws.Range("E17").Value = "  -2.51%  "
$ws.Range("D18").Value = "1.982.24"
$ws.Range("E18").Value = "  -3.56%  "
$ws.Range("D19").Value = "36.354.23"
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").Value = "70.65"
$ws.Range("E20").Value = "  -4.26%  "
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "5.31"
$ws.Range("E22").Value = "  -2.71%  "
$ws.Range("D23").Value = "234.37"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -3.21%  "
$ws.Range("E26").Value = "  -3.87%  "
$ws.Range("D27").Value = "10.15"
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("D28").Value = "161.87"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "19.83"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("D30").Value = "0.129"
$ws.Range("E30").Value = "  +8.84%  "
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").Value = "4.91"
$ws.Range("E33").Value = "  -6.94%  "
$ws.Range("D34").Value = "0.0631"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("E35").Value = "  -6.11%  "
$ws.Range("D36").Value = "6.30"
$ws.Range("E36").Value = "  +4.31%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  -7.82%  "
$ws.Range("E39").Value = "  -4.54%  "
$ws.Range("E40").Value = "  +3.55%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "0.0960"
$ws.Range("E42").Value = "  -7.59%  "
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("E45").Value = "  -4.96%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "16.24"
$ws.Range("E46").Value = "  -5.52%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "92.58"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("E48").Value = "  -5.40%  "
$ws.Range("D49").Value = "1.375.53"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("E50").Value = "  -3.19%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.95"
$ws.Range("E51").Value = "  +9.49%  "
